$wb = $excel.ActiveWorkbook

# --- Rename the second sheet ---
$ws2 = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$ws2.Name = "Include #0"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Date value (row 8)
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Contact value (row 9 header is row 10 actually; Contact label is A10)
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row for "Jurisdiction" right after the Contact row (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Match the border/wrap/alignment formatting used by the rest of the data rows
# (reuse the existing thin grey border + top-aligned wrapped text look).
$newRow = $ws.Range("A11:B11")
$newRow.Borders.Color = 8421504
$newRow.Borders.LineStyle = 1
$newRow.VerticalAlignment = -4160
$newRow.WrapText = $true

$ws.Range("A11").Value = "Jurisdiction"

Write-Output "done"
